# Add the new "Max Consecutive Ones" (LeetCode 485) entry as row 12 to the
# "数组" (Arrays) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("数组")

# No. / leetcode id
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 485

# Shared-string write order matters for de-dup index assignment: the
# solution-steps cell (D12) must be written before the problem-statement
# cell (C12) so the new shared strings land at indices 182/183/184 in the
# same order the source workbook has them.
$ws.Range("D12").Value = "1 初始化max=0，累加变量count=0，迭代数组元素`n2 如果元素是1，就count++`n3 如果元素是0，count=0`n4 无论是否是1，取max与count的最大值赋值给max，这个办法比较耗时`n优化：只有是0的时候取max，迭代结束再次判断max`n5 迭代结束，max就是最终结果"

$ws.Range("C12").Value = "给定一个二进制数组， 计算其中最大连续1的个数。 `n 示例 1: `n输入: [1,1,0,1,1,1]`n输出: 3`n解释: 开头的两位和最后的三位都是连续1，所以最大连续1的个数是 3.`n 注意： `n 输入的数组只包含 0 和1。 `n 输入数组的长度是正整数，且不超过 10,000。 `n Related Topics 数组"

$ws.Range("E12").Value = "保留上次最大值`n计数器清零"

$ws.Range("F12").Value = "O(N)"
$ws.Range("G12").Value = "O(1)"

# Match Excel's auto-fit wrapped-text row height for the new row.
$ws.Rows.Item(12).RowHeight = 286

# Scroll the sheet to reveal the new row, mirroring the author's view state.
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("D16").Select()
